$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.857.21'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.938.29'
$ws.Range("E3").Value = '  -0.98%  '
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.66'
$ws.Range("E5").Value = '  -1.03%  '
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4913'
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2939'
$ws.Range("E8").Value = '  -1.83%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06886'
$ws.Range("E9").Value = '  +0.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.24'
$ws.Range("E10").Value = '  +0.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '104.96'
$ws.Range("E11").Value = '  -3.21%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.944.94'
$ws.Range("E12").Value = '  -0.49%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07774'
$ws.Range("E13").Value = '  +0.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.357'
$ws.Range("E14").Value = '  -1.97%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.7011'
$ws.Range("E15").Value = '  -1.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '274.10'
$ws.Range("E16").Value = '  -3.52%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.859.58'
$ws.Range("E17").Value = '  -0.95%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007723'
$ws.Range("E18").Value = '  -0.69%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.10'
$ws.Range("E19").Value = '  -1.29%  '
$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  -0.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.579'
$ws.Range("E21").Value = '  +1.38%  '
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.193.04'
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  -0.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.531'
$ws.Range("E24").Value = '  -0.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.857'
$ws.Range("E25").Value = '  +0.28%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.18'
$ws.Range("E26").Value = '  -2.18%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.61'
$ws.Range("E27").Value = '  -2.50%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.163'
$ws.Range("E28").Value = '  -2.65%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1044'
$ws.Range("E29").Value = '  -0.68%  '
$ws.Range("E30").Value = '  -2.76%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.557'
$ws.Range("E31").Value = '  -1.82%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.556'
$ws.Range("E32").Value = '  -0.63%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.377'
$ws.Range("E33").Value = '  -1.80%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04894'
$ws.Range("E34").Value = '  -1.70%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7595'
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.151'
$ws.Range("E36").Value = '  -2.48%  '
$ws.Range("E37").Value = '  -0.24%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.736'
$ws.Range("E38").Value = '  -0.12%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02008'
$ws.Range("E39").Value = '  -1.81%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.658'
$ws.Range("E40").Value = '  -2.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.548'
$ws.Range("E41").Value = '  +1.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '77.96'
$ws.Range("E42").Value = '  +6.41%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.093'
$ws.Range("E43").Value = '  -3.71%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9132'
$ws.Range("E44").Value = '  +3.29%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4441'
$ws.Range("E45").Value = '  -1.52%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '107.83'
$ws.Range("E46").Value = '  -1.52%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9991'
$ws.Range("E47").Value = '  -0.36%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.723'
$ws.Range("E48").Value = '  -6.42%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '995.45'
$ws.Range("E49").Value = '  +3.53%  '
$ws.Range("E50").Value = '  -1.52%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '36.14'
$ws.Range("E51").Value = '  +1.17%  '
